# Apply cryptos list refresh (coin prices / 1h volume %) as captured in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'63.236.20"
$ws.Cells.Item(2, 5).Value = "  -1.31%  "

$ws.Cells.Item(3, 4).Value = "'3.046.50"
$ws.Cells.Item(3, 5).Value = "  -3.33%  "

$ws.Cells.Item(4, 5).Value = "  +0.05%  "

$ws.Cells.Item(5, 4).Value = "'560.59"
$ws.Cells.Item(5, 5).Value = "  -1.25%  "

$ws.Cells.Item(6, 4).Value = "'157.80"
$ws.Cells.Item(6, 5).Value = "  -5.07%  "

$ws.Cells.Item(7, 5).Value = "  +0.16%  "

$ws.Cells.Item(8, 4).Value = "'0.560"
$ws.Cells.Item(8, 5).Value = "  -6.26%  "

$ws.Cells.Item(9, 4).Value = "'3.051.43"
$ws.Cells.Item(9, 5).Value = "  -3.17%  "

$ws.Cells.Item(10, 4).Value = "'0.114"
$ws.Cells.Item(10, 5).Value = "  -3.08%  "

$ws.Cells.Item(11, 4).Value = "'6.50"

$ws.Cells.Item(12, 4).Value = "'0.370"
$ws.Cells.Item(12, 5).Value = "  -3.47%  "

$ws.Cells.Item(13, 4).Value = "'3.600.31"
$ws.Cells.Item(13, 5).Value = "  -2.64%  "

$ws.Cells.Item(14, 5).Value = "  -2.29%  "

$ws.Cells.Item(15, 4).Value = "'63.435.27"
$ws.Cells.Item(15, 5).Value = "  -1.14%  "

$ws.Cells.Item(16, 4).Value = "'23.98"
$ws.Cells.Item(16, 5).Value = "  -4.22%  "

$ws.Cells.Item(17, 4).Value = "'3.074.09"
$ws.Cells.Item(17, 5).Value = "  -2.65%  "

$ws.Cells.Item(18, 4).Value = "'0.0000151"
$ws.Cells.Item(18, 5).Value = "  -2.66%  "

$ws.Cells.Item(19, 4).Value = "'397.35"
$ws.Cells.Item(19, 5).Value = "  -3.78%  "

$ws.Cells.Item(20, 4).Value = "'5.10"
$ws.Cells.Item(20, 5).Value = "  -2.39%  "

$ws.Cells.Item(21, 4).Value = "'12.07"

$ws.Cells.Item(22, 4).Value = "'6.79"
$ws.Cells.Item(22, 5).Value = "  -4.11%  "

$ws.Cells.Item(23, 4).Value = "'0.998"
$ws.Cells.Item(23, 5).Value = "  +0.00%  "

$ws.Cells.Item(24, 4).Value = "'66.59"
$ws.Cells.Item(24, 5).Value = "  -4.04%  "

$ws.Cells.Item(25, 4).Value = "'0.469"
$ws.Cells.Item(25, 5).Value = "  -4.97%  "

$ws.Cells.Item(26, 4).Value = "'0.189"
$ws.Cells.Item(26, 5).Value = "  -7.15%  "

$ws.Cells.Item(27, 4).Value = "'0.0₃0984"
$ws.Cells.Item(27, 5).Value = "  -2.31%  "

$ws.Cells.Item(28, 4).Value = "'8.82"
$ws.Cells.Item(28, 5).Value = "  +1.27%  "

$ws.Cells.Item(29, 5).Value = "  +0.17%  "

$ws.Cells.Item(30, 5).Value = "  +0.11%  "

$ws.Cells.Item(31, 4).Value = "'1.77"
$ws.Cells.Item(31, 5).Value = "  -2.42%  "

$ws.Cells.Item(32, 4).Value = "'20.72"
$ws.Cells.Item(32, 5).Value = "  -3.68%  "

$ws.Cells.Item(33, 4).Value = "'162.53"
$ws.Cells.Item(33, 5).Value = "  +4.58%  "

$ws.Cells.Item(34, 4).Value = "'4.74"
$ws.Cells.Item(34, 5).Value = "  -4.39%  "

$ws.Cells.Item(35, 4).Value = "'6.06"
$ws.Cells.Item(35, 5).Value = "  -3.46%  "

$ws.Cells.Item(36, 4).Value = "'1.09"
$ws.Cells.Item(36, 5).Value = "  -1.13%  "

$ws.Cells.Item(37, 4).Value = "'1.32"
$ws.Cells.Item(37, 5).Value = "  -2.01%  "

$ws.Cells.Item(38, 4).Value = "'1.62"
$ws.Cells.Item(38, 5).Value = "  -3.34%  "

$ws.Cells.Item(39, 4).Value = "'2.532.73"
$ws.Cells.Item(39, 5).Value = "  -5.50%  "

$ws.Cells.Item(40, 2).Value = "Filecoin"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(40, 4).Value = "'3.99"
$ws.Cells.Item(40, 5).Value = "  -3.86%  "

$ws.Cells.Item(41, 2).Value = "EnergySwap"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(41, 4).Value = "'22.97"
$ws.Cells.Item(41, 5).Value = "  -3.18%  "

$ws.Cells.Item(42, 4).Value = "'37.73"
$ws.Cells.Item(42, 5).Value = "  -3.00%  "

$ws.Cells.Item(43, 4).Value = "'0.671"
$ws.Cells.Item(43, 5).Value = "  -5.50%  "

$ws.Cells.Item(44, 4).Value = "'0.0601"
$ws.Cells.Item(44, 5).Value = "  -1.66%  "

$ws.Cells.Item(45, 4).Value = "'0.0248"
$ws.Cells.Item(45, 5).Value = "  -3.94%  "

$ws.Cells.Item(46, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Cells.Item(46, 4).Value = "'1.00"
$ws.Cells.Item(46, 5).Value = "  +0.00%  "

$ws.Cells.Item(47, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(47, 4).Value = "'20.36"
$ws.Cells.Item(47, 5).Value = "  -3.03%  "

$ws.Cells.Item(48, 2).Value = "Bittensor"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Cells.Item(48, 4).Value = "'275.99"
$ws.Cells.Item(48, 5).Value = "  -3.75%  "

$ws.Cells.Item(49, 2).Value = "RenderToken"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(49, 4).Value = "'4.96"
$ws.Cells.Item(49, 5).Value = "  -6.73%  "

$ws.Cells.Item(50, 4).Value = "'0.0954"
$ws.Cells.Item(50, 5).Value = "  -2.97%  "

$ws.Cells.Item(51, 5).Value = "  +0.32%  "
